$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    3 = @{ B = "15-24"; C = "off";   D = "off";   E = "off";   F = "7-16";  G = "10-19" }
    4 = @{ B = "15-24"; C = "15-24"; D = "off";   E = "10-19"; F = "off";   G = "7-16"  }
    5 = @{ B = "7-16";  C = "off";   D = "off";   E = "15-24"; F = "10-19"; G = "15-24" }
    6 = @{ B = "10-19"; C = "15-24"; D = "15-24"; E = "off";   F = "off";   G = "7-16"  }
    7 = @{ B = "15-24"; C = "off";   D = "15-24"; E = "7-16";  F = "10-19"; G = "off"   }
    8 = @{ B = "15-24"; C = "off";   D = "10-19"; E = "off";   F = "15-24"; G = "7-16"  }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
